{"js": "// The letter originally reads:\n//   \"I am the owner of Perfect Paws L.L.C., and I would be delighted for\n//    you to assist me in setting up a computerized database billing system.\"\n// The edit swaps the company name \"Perfect Paws L.L.C.,\" for the generic\n// phrase \"your local pet company\", keeping the rest of the sentence intact.\nconst body = context.document.body;\n\nconst results = body.search(\"Perfect Paws L.L.C.,\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text 'Perfect Paws L.L.C.,' in document body.\");\n}\n\n// Replace just that run of text; Word re-splits the surrounding runs and\n// keeps their original character formatting (Arial / black) automatically.\nresults.items[0].insertText(\"your local pet company\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The letter originally reads:\n#   \"I am the owner of Perfect Paws L.L.C., and I would be delighted for\n#    you to assist me in setting up a computerized database billing system.\"\n# The edit swaps the company name \"Perfect Paws L.L.C.,\" for the generic\n# phrase \"your local pet company\", keeping the rest of the sentence intact.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Perfect Paws L.L.C.,\"\n$find.Replacement.Text = \"your local pet company\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)\n"}
